# Apply edit: update DigiKey part number text, add column G with formula F*2+1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the shared-string text for E13 (DigiKey part number for U4)
$ws.Range("E13").Value = "NCP1117ST33T3GOSCT-ND"

# 2. Add new column G with formula =F{row}*2+1 for rows 1-13
$ws.Range("G1").Formula = "=F1*2+1"
$ws.Range("G2:G13").Formula = "=F2*2+1"

# 3. Update selection to the full used range A1:G13
$ws.UsedRange.Select()
